# Auto stash before merge of "master" and "origin/master"
#
# Replace the single "DogZone" sample customer row with real customer
# records, and append three more customer rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Company Name | B Customer ID | C Address | D Primary Title
#          E Primary Contact | F Pri. Contact Number | G Primary Email Address
#          H Secondary Title | I Secondary Contact | J Sec Contact Number
#          K Secondary Email Address

# Customer ID / phone-number style columns look like numbers, but in the
# source workbook they are plain shared-string text. Force those specific
# cells to Text before writing so Excel doesn't reinterpret them as
# numeric values, then drop the temporary formatting again so the cells
# end up unstyled, exactly like every other text cell in the sheet.
$textCells = @("B2", "F2", "B3", "B4", "B5", "F5")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - True Yoga Pte Ltd (replaces the old "DogZone" placeholder row)
# D2/H2/I2/J2/K2 are already blank in the source sheet, so they are left
# untouched here rather than re-assigned.
$ws.Range("A2").Value = "True Yoga Pte Ltd"
$ws.Range("B2").Value = "1300001"
$ws.Range("C2").Value = "No.8 Claymore Hill, #02-03, 8 on Claymore Singapore 229572"
$ws.Range("E2").Value = "Davis Lee"
$ws.Range("F2").Value = "66727236"
$ws.Range("G2").Value = "davis.lee@trueyoga.com.sg"

# Row 3 - Gizmo Square
$ws.Range("A3").Value = "Gizmo Square"
$ws.Range("B3").Value = "1300002"
$ws.Range("E3").Value = "Curtis Chen"
$ws.Range("G3").Value = "curtis@gizmosquare.com.sg"
$ws.Range("I3").Value = "Yolo"
$ws.Range("K3").Value = "kenneth@gizmosquare.com.sg"

# Row 4 - Hi / Bye
$ws.Range("A4").Value = "Hi"
$ws.Range("B4").Value = "1300003"
$ws.Range("C4").Value = "Bye"

# Row 5 - Bye / Givemee
$ws.Range("A5").Value = "Bye"
$ws.Range("B5").Value = "1300004"
$ws.Range("C5").Value = "Givemee"
$ws.Range("E5").Value = "hi"
$ws.Range("F5").Value = "98123423"
$ws.Range("G5").Value = "Panda@hotmail.com"

# Strip the temporary text formatting back off so the cells don't retain
# an explicit style index (matches the original workbook where every
# cell is unstyled).
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
